$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column B (Apellido) - shifts old B->C, old C->D
$ws.Columns("B:B").Insert()

# Header row
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Apellido"

# Surnames for existing rows
$ws.Range("B3").Value = "Alvarez"
$ws.Range("B5").Value = "Cienfuegos"
$ws.Range("B6").Value = "Dominguez"

# Update Borja's NIF value (now in D4) from number to text
$ws.Range("D4").Value = "22222222B"

$ws.Range("E6").Select() | Out-Null
